$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.338.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "'1.830.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'329.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.4446"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").Value = "'0.3772"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "'44.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "'0.07749"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").Value = "'1.135"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'22.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "'1.001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'6.364"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "'7.568"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'1.837.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "'93.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +15.49%  "
$ws.Range("D18").Value = "'0.00001086"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'0.06459"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "'17.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "'6.363"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").Value = "'0.5409"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'28.391.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "'11.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Value = "'2.245"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.05%  "
$ws.Range("D27").Value = "'20.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").Value = "'155.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").Value = "'2.374"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "'2.043.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").Value = "'128.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").Value = "'1.213"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.69%  "
$ws.Range("D33").Value = "'5.950"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("D34").Value = "'0.09307"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "'3.688"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.28%  "
$ws.Range("D36").Value = "'13.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.22%  "
$ws.Range("D37").Value = "'0.02358"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").Value = "'0.2198"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").Value = "'5.220"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'0.6627"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "'0.06256"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'8.214"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("D43").Value = "'1.200"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "'0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45 and 46: coin names swapped (WEMIXTOKEN <-> EnergySwap) along with updated price/volume values
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'14.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.64%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.393"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.05%  "

$ws.Range("D47").Value = "'0.6138"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "'3.787"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'2.063"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").Value = "'127.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "'0.07017"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
